# Applies the "remove the Daum row" edit:
#  - B8 changes from 2.5 to 3.5
#  - Row 10 (Daum #다음 entry) is deleted entirely, which also drops the
#    now-unused shared strings and its hyperlink/relationship
#  - The active selection moves from E7 to B7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the current per-cell formatting of the hyperlinked D2:D9 cells so it
# can be restored later - removing/recreating hyperlinks resets cell style.
$ws.Range("D2:D9").Copy()
$ws.Range("H2:H9").PasteSpecial(-4122) # xlPasteFormats

# Update the interval value for the "normal" row (row 8)
$ws.Range("B8").Value = 3.5

# Excel recreates hyperlink relationships on save; clear the existing set so
# the row-10 (Daum) hyperlink does not survive as a dangling reference once
# its row is removed below.
$ws.Hyperlinks.Delete()

# Remove the whole 10th row (A10:E10, the "#다음" / daum.net entry)
$ws.Rows("10:10").Delete()

# Recreate the hyperlinks for the rows that remain, same order/targets as
# the original workbook
$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.asiatime.co.kr/search?searchText=%EC%95%84%EC%8B%9C%EC%95%84%ED%83%80%EC%9E%84%EC%A6%88")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.segyebiz.com/newsList/0500000000000")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.ichannela.com/news/main/news_part.do?catecode=000400")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://news.mtn.co.kr/category-news/all")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://www.yakup.com/news/index.html")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://m.dnews.co.kr/m_home/index.html")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://m.inews24.com/l/recency")
$ws.Hyperlinks.Add($ws.Range("D9"), "https://www.hanaw.com/main/research/trends/RC_060600_P1.cmd")

# Restore the original cell formatting (Hyperlinks.Add applies its own style)
$ws.Range("H2:H9").Copy()
$ws.Range("D2:D9").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H2:H9").Clear()

# Move the active selection from E7 to B7
$null = $ws.Range("B7").Select()
